$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the awareness_1 question text in cell A4: remove ", która pojawi się
# poniżej" after "klawiatury" per Balas's comments.
$newText = "Czy zauważyłeś coś charakterstycznego w słowach i obrazach, które pojawiały się z określonymi postaciami z kreskówek?`n`nWpisz swoją odpowiedź za pomocą klawiatury.`n`nNaciśnij Enter by przejść do następnego pytania.`n`nTwoja odpowiedź musi mieć minimum 20 znaków."
$ws.Range("A4").Value = $newText

# Remove sheet protection that was previously applied.
$ws.Unprotect()
